# Auto-generated edit script applying numeric corrections to leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1232.8695
$ws.Range("I19").Value = 827.8570999999999
$ws.Range("J19").Value = 1410.0625
$ws.Range("K19").Value = 827.8570999999999
$ws.Range("L19").Value = 1410.0625
$ws.Range("M19").Value = -652.8570999999999
$ws.Range("N19").Value = -1760.0625
$ws.Range("H107").Value = 758.86957
$ws.Range("I107").Value = 470.8889
$ws.Range("K107").Value = 470.8889
$ws.Range("M107").Value = 1449.1111
$ws.Range("H132").Value = 15923066
$ws.Range("I132").Value = 22224146
$ws.Range("J132").Value = 170369.17
$ws.Range("K132").Value = 66672438
$ws.Range("L132").Value = 511107.51
$ws.Range("M132").Value = -66669908
$ws.Range("N132").Value = -516167.51
$ws.Range("H137").Value = 6868.3687
$ws.Range("I137").Value = 1589.4546
$ws.Range("J137").Value = 14126.875
$ws.Range("K137").Value = 4768.3638
$ws.Range("L137").Value = 42380.625
$ws.Range("M137").Value = -2218.3638
$ws.Range("N137").Value = -47480.625
$ws.Range("H138").Value = 2693.2104
$ws.Range("J138").Value = 3643.653
$ws.Range("L138").Value = 10930.959
$ws.Range("N138").Value = -21210.959
$ws.Range("H139").Value = 289999.75
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("H141").Value = 1397.8572
$ws.Range("I141").Value = 1397.8572
$ws.Range("K141").Value = 4193.571599999999
$ws.Range("M141").Value = 986.4284000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 34485980
$ws.Range("I2").Value = 50003070
$ws.Range("K2").Value = 50003070
$ws.Range("M2").Value = -50002957
$ws.Range("H10").Value = 5777
$ws.Range("I10").Value = 5777
$ws.Range("K10").Value = 5777
$ws.Range("M10").Value = -5607
$ws.Range("H110").Value = 4889.6216
$ws.Range("I110").Value = 4617.931
$ws.Range("K110").Value = 4617.931
$ws.Range("M110").Value = -2572.931
$ws.Range("H116").Value = 34485980
$ws.Range("I116").Value = 50003070
$ws.Range("K116").Value = 50003070
$ws.Range("M116").Value = -50000776
$ws.Range("H122").Value = 2091.5715
$ws.Range("I122").Value = 1940.1666
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5820.4998
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3370.4998
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 4899.5835
$ws.Range("I132").Value = 5009.778
$ws.Range("K132").Value = 15029.334
$ws.Range("M132").Value = -12499.334
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 34485980
$ws.Range("I3").Value = 50003070
$ws.Range("K3").Value = 50003070
$ws.Range("M3").Value = -50002956
$ws.Range("H134").Value = 1695.9517
$ws.Range("I134").Value = 1646.6271
$ws.Range("K134").Value = 4939.8813
$ws.Range("M134").Value = -2404.8813
$ws.Range("H140").Value = 103999.336
$ws.Range("J140").Value = 103999.336
$ws.Range("L140").Value = 103999.336
$ws.Range("N140").Value = -114359.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1020.7273
$ws.Range("I16").Value = 1005
$ws.Range("J16").Value = 1048.25
$ws.Range("K16").Value = 1005
$ws.Range("L16").Value = 1048.25
$ws.Range("M16").Value = -718
$ws.Range("N16").Value = -1622.25
$ws.Range("H31").Value = 22681.666
$ws.Range("I31").Value = 39784.406
$ws.Range("K31").Value = 39784.406
$ws.Range("M31").Value = -39489.406
$ws.Range("H34").Value = 22681.666
$ws.Range("I34").Value = 39784.406
$ws.Range("K34").Value = 39784.406
$ws.Range("M34").Value = -39582.406
$ws.Range("H58").Value = 3123.5454
$ws.Range("I58").Value = 2773
$ws.Range("K58").Value = 2773
$ws.Range("M58").Value = -2570
$ws.Range("H105").Value = 1991.1666
$ws.Range("I105").Value = 1511.5454
$ws.Range("K105").Value = 1511.5454
$ws.Range("M105").Value = 235.4546
$ws.Range("H107").Value = 1040.1282
$ws.Range("I107").Value = 781.86365
$ws.Range("J107").Value = 1374.3529
$ws.Range("K107").Value = 781.86365
$ws.Range("L107").Value = 1374.3529
$ws.Range("M107").Value = 1138.13635
$ws.Range("N107").Value = -5214.3529
$ws.Range("H113").Value = 1020.7273
$ws.Range("I113").Value = 1005
$ws.Range("J113").Value = 1048.25
$ws.Range("K113").Value = 1005
$ws.Range("L113").Value = 1048.25
$ws.Range("M113").Value = 1165
$ws.Range("N113").Value = -5388.25
$ws.Range("H115").Value = 16500
$ws.Range("J115").Value = 16500
$ws.Range("L115").Value = 16500
$ws.Range("N115").Value = -18850
$ws.Range("H132").Value = 138188.67
$ws.Range("I132").Value = 252548.16
$ws.Range("J132").Value = 3648.1177
$ws.Range("K132").Value = 757644.48
$ws.Range("L132").Value = 10944.3531
$ws.Range("M132").Value = -755114.48
$ws.Range("N132").Value = -16004.3531
$ws.Range("H134").Value = 22377.365
$ws.Range("I134").Value = 18396.879
$ws.Range("J134").Value = 37213.727
$ws.Range("K134").Value = 55190.637
$ws.Range("L134").Value = 111641.181
$ws.Range("M134").Value = -52655.637
$ws.Range("N134").Value = -116711.181
$ws.Range("H136").Value = 3123.5454
$ws.Range("I136").Value = 2773
$ws.Range("K136").Value = 8319
$ws.Range("M136").Value = -5769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 733.6667
$ws.Range("I22").Value = 201
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 603
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -434
$ws.Range("N22").Value = -3338
$ws.Range("H27").Value = 733.6667
$ws.Range("I27").Value = 201
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 603
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -501
$ws.Range("N27").Value = -3204
$ws.Range("H51").Value = 430.4
$ws.Range("I51").Value = 433
$ws.Range("J51").Value = 426.5
$ws.Range("K51").Value = 1299
$ws.Range("L51").Value = 1279.5
$ws.Range("M51").Value = -839
$ws.Range("N51").Value = -2199.5
$ws.Range("H132").Value = 1155.9487
$ws.Range("I132").Value = 1090.7354
$ws.Range("K132").Value = 9816.6186
$ws.Range("M132").Value = -7286.6186
$ws.Range("H141").Value = 384713.6
$ws.Range("I141").Value = 230883.75
$ws.Range("K141").Value = 692651.25
$ws.Range("M141").Value = -687471.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10272.259
$ws.Range("I70").Value = 10512.25
$ws.Range("J70").Value = 10171.211
$ws.Range("K70").Value = 10512.25
$ws.Range("L70").Value = 10171.211
$ws.Range("M70").Value = -10242.25
$ws.Range("N70").Value = -10711.211
$ws.Range("H73").Value = 10272.259
$ws.Range("I73").Value = 10512.25
$ws.Range("J73").Value = 10171.211
$ws.Range("K73").Value = 10512.25
$ws.Range("L73").Value = 10171.211
$ws.Range("M73").Value = -9576.25
$ws.Range("N73").Value = -12043.211
$ws.Range("H80").Value = 6098.5386
$ws.Range("J80").Value = 8384
$ws.Range("L80").Value = 8384
$ws.Range("N80").Value = -10380
$ws.Range("H83").Value = 6098.5386
$ws.Range("J83").Value = 8384
$ws.Range("L83").Value = 41920
$ws.Range("N83").Value = -51904
$ws.Range("H122").Value = 2075.0715
$ws.Range("I122").Value = 1007.1111
$ws.Range("J122").Value = 3997.4
$ws.Range("K122").Value = 3021.3333
$ws.Range("L122").Value = 11992.2
$ws.Range("M122").Value = -571.3332999999998
$ws.Range("N122").Value = -16892.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11412.818
$ws.Range("I22").Value = 1660.2
$ws.Range("J22").Value = 19540
$ws.Range("K22").Value = 1660.2
$ws.Range("L22").Value = 19540
$ws.Range("M22").Value = -1365.2
$ws.Range("N22").Value = -20130
$ws.Range("H27").Value = 11412.818
$ws.Range("I27").Value = 1660.2
$ws.Range("J27").Value = 19540
$ws.Range("K27").Value = 1660.2
$ws.Range("L27").Value = 19540
$ws.Range("M27").Value = -1553.2
$ws.Range("N27").Value = -19754
$ws.Range("H46").Value = 1120.6666
$ws.Range("I46").Value = 1105
$ws.Range("J46").Value = 1199
$ws.Range("K46").Value = 1105
$ws.Range("L46").Value = 1199
$ws.Range("M46").Value = -917
$ws.Range("N46").Value = -1575
$ws.Range("H82").Value = 2666.5789
$ws.Range("I82").Value = 2862.2
$ws.Range("J82").Value = 2449.2222
$ws.Range("K82").Value = 2862.2
$ws.Range("L82").Value = 2449.2222
$ws.Range("M82").Value = -2501.2
$ws.Range("N82").Value = -3171.2222
$ws.Range("H85").Value = 2666.5789
$ws.Range("I85").Value = 2862.2
$ws.Range("J85").Value = 2449.2222
$ws.Range("K85").Value = 2862.2
$ws.Range("L85").Value = 2449.2222
$ws.Range("M85").Value = -1614.2
$ws.Range("N85").Value = -4945.2222
$ws.Range("H105").Value = 41633.332
$ws.Range("I105").Value = 29900
$ws.Range("J105").Value = 47500
$ws.Range("K105").Value = 29900
$ws.Range("L105").Value = 47500
$ws.Range("M105").Value = -26406
$ws.Range("N105").Value = -54488
$ws.Range("H132").Value = 2591.149
$ws.Range("I132").Value = 2367.2
$ws.Range("J132").Value = 3870.8572
$ws.Range("K132").Value = 7101.599999999999
$ws.Range("L132").Value = 11612.5716
$ws.Range("M132").Value = -4571.599999999999
$ws.Range("N132").Value = -16672.5716
$ws.Range("H139").Value = 99996.664
$ws.Range("J139").Value = 99996.664
$ws.Range("L139").Value = 99996.664
$ws.Range("N139").Value = -110276.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 12001.333
$ws.Range("I9").Value = 14252
$ws.Range("J9").Value = 7500
$ws.Range("K9").Value = 14252
$ws.Range("L9").Value = 7500
$ws.Range("M9").Value = -14112
$ws.Range("N9").Value = -7780
$ws.Range("H113").Value = 974.36365
$ws.Range("I113").Value = 968.6667
$ws.Range("K113").Value = 2906.0001
$ws.Range("M113").Value = -736.0001000000002
$ws.Range("H122").Value = 2084.875
$ws.Range("J122").Value = 2363.1
$ws.Range("L122").Value = 7089.299999999999
$ws.Range("N122").Value = -11989.3

